$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet: new source file uuid + refreshed "Latest HO Xliff
# Generate Date" timestamp.
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "13a65221-79e1-4eff-a0b3-20040817c627.md"
$wsOverview.Range("B2").Value = "e2e\13a65221-79e1-4eff-a0b3-20040817c627.md"
$wsOverview.Range("G2").Value = "2016-08-25 17:02:38"

foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\13a65221-79e1-4eff-a0b3-20040817c627.md"
}

# ---------------------------------------------------------------------------
# zh-cn handoff-status sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "13a65221-79e1-4eff-a0b3-20040817c627.md"
$wsZhCn.Range("G2").Value = "13a65221-79e1-4eff-a0b3-20040817c627.289ac71cec9544d3a1e7877bd6d7061c5f1039f0.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-25 17:02:34"

# Target/handback file is no longer applicable for this handoff round - the
# "Latest Target File" hyperlink+value is cleared (style reset to Normal)
# and the "Latest Handback File" value is cleared, while only the first
# (source-file) hyperlink on the row is kept - with a refreshed display name.
foreach ($h in $wsZhCn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$I`$2") {
        $h.Delete()
    } else {
        $h.TextToDisplay = "13a65221-79e1-4eff-a0b3-20040817c627.md"
    }
}
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""

# No handback has happened yet, so the handback datetime resets to the
# default (unset) datetime.
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Column widths for "Latest Target File" / "Latest Handback File" shrink
# now that they hold short/empty values instead of full file names.
$wsZhCn.Columns.Item(9).ColumnWidth = 17.8
$wsZhCn.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------------
# de-de handoff-status sheet (mirrors zh-cn)
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "13a65221-79e1-4eff-a0b3-20040817c627.md"
$wsDeDe.Range("G2").Value = "13a65221-79e1-4eff-a0b3-20040817c627.289ac71cec9544d3a1e7877bd6d7061c5f1039f0.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-25 17:02:38"

foreach ($h in $wsDeDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$I`$2") {
        $h.Delete()
    } else {
        $h.TextToDisplay = "13a65221-79e1-4eff-a0b3-20040817c627.md"
    }
}
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""

$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Columns.Item(9).ColumnWidth = 17.8
$wsDeDe.Columns.Item(10).ColumnWidth = 20.8
